$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "93.231.45"
Set-TextValue $ws "E2" "  +1.75%  "

Set-TextValue $ws "D3" "3.098.97"
Set-TextValue $ws "E3" "  -0.59%  "

Set-TextValue $ws "E4" "  +0.00%  "

Set-TextValue $ws "D5" "236.50"
Set-TextValue $ws "E5" "  -3.91%  "

Set-TextValue $ws "D6" "612.87"
Set-TextValue $ws "E6" "  -0.68%  "

Set-TextValue $ws "D7" "1.12"
Set-TextValue $ws "E7" "  +1.90%  "

Set-TextValue $ws "D8" "0.388"
Set-TextValue $ws "E8" "  +1.33%  "

Set-TextValue $ws "E9" "  -0.04%  "

Set-TextValue $ws "D10" "0.828"
Set-TextValue $ws "E10" "  +13.02%  "

Set-TextValue $ws "D11" "3.094.44"
Set-TextValue $ws "E11" "  -0.73%  "

Set-TextValue $ws "E12" "  -2.98%  "

Set-TextValue $ws "E13" "  -2.87%  "

Set-TextValue $ws "E14" "  +0.37%  "

Set-TextValue $ws "D15" "93.005.24"
Set-TextValue $ws "E15" "  +1.57%  "

Set-TextValue $ws "D16" "5.41"
Set-TextValue $ws "E16" "  -3.21%  "

Set-TextValue $ws "D17" "3.674.66"
Set-TextValue $ws "E17" "  -0.66%  "

Set-TextValue $ws "D18" "3.104.06"
Set-TextValue $ws "E18" "  +0.65%  "

Set-TextValue $ws "D19" "3.68"
Set-TextValue $ws "E19" "  -0.29%  "

Set-TextValue $ws "D20" "14.61"
Set-TextValue $ws "E20" "  -1.59%  "

Set-TextValue $ws "D21" "5.95"
Set-TextValue $ws "E21" "  +2.67%  "

Set-TextValue $ws "D22" "441.11"
Set-TextValue $ws "E22" "  -1.16%  "

Set-TextValue $ws "E23" "  -1.68%  "

Set-TextValue $ws "D24" "9.03"
Set-TextValue $ws "E24" "  -5.05%  "

Set-TextValue $ws "D25" "8.23"
Set-TextValue $ws "E25" "  +4.66%  "

Set-TextValue $ws "D26" "5.68"
Set-TextValue $ws "E26" "  -2.87%  "

Set-TextValue $ws "D27" "12.72"
Set-TextValue $ws "E27" "  +8.36%  "

Set-TextValue $ws "E28" "  -2.62%  "

Set-TextValue $ws "E29" "  -0.22%  "

Set-TextValue $ws "D30" "0.247"
Set-TextValue $ws "E30" "  +4.96%  "

Set-TextValue $ws "D31" "0.182"
Set-TextValue $ws "E31" "  +8.50%  "

Set-TextValue $ws "E32" "  -13.86%  "

Set-TextValue $ws "D33" "9.17"
Set-TextValue $ws "E33" "  -1.35%  "

Set-TextValue $ws "E34" "  +0.71%  "

Set-TextValue $ws "D35" "7.91"
Set-TextValue $ws "E35" "  +0.72%  "

Set-TextValue $ws "E36" "  -9.98%  "

Set-TextValue $ws "D37" "25.88"
Set-TextValue $ws "E37" "  -1.23%  "

Set-TextValue $ws "E38" "  -4.35%  "

Set-TextValue $ws "E39" "  -2.16%  "

Set-TextValue $ws "B40" "PolygonEcosystemToken"
Set-TextValue $ws "C40" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D40" "0.444"
Set-TextValue $ws "E40" "  +1.15%  "

Set-TextValue $ws "B41" "WhiteBITCoin"
Set-TextValue $ws "C41" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D41" "23.97"
Set-TextValue $ws "E41" "  +8.04%  "

Set-TextValue $ws "B42" "Bittensor"
Set-TextValue $ws "C42" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D42" "475.99"
Set-TextValue $ws "E42" "  -2.94%  "

Set-TextValue $ws "B43" "Fetch.AI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D43" "1.29"
Set-TextValue $ws "E43" "  -1.19%  "

Set-TextValue $ws "D44" "3.27"
Set-TextValue $ws "E44" "  -3.96%  "

Set-TextValue $ws "E45" "  +0.04%  "

Set-TextValue $ws "D46" "158.97"
Set-TextValue $ws "E46" "  +0.79%  "

Set-TextValue $ws "E47" "  -1.32%  "

Set-TextValue $ws "E48" "  -2.45%  "

Set-TextValue $ws "E49" "  -2.21%  "

Set-TextValue $ws "D50" "43.82"
Set-TextValue $ws "E50" "  -0.52%  "

Set-TextValue $ws "D51" "4.37"
Set-TextValue $ws "E51" "  -0.32%  "
